$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data table. Insert a fresh
# row at row 348 (pushing the existing row 348 and everything below it
# down by one row, so the previous last row 394 becomes row 395) and
# populate it with the new record's values.
$ws.Rows.Item(348).Insert()

$ws.Cells.Item(348, 1).Value  = 9
$ws.Cells.Item(348, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(348, 3).Value  = "Metropolitana"
$ws.Cells.Item(348, 4).Value  = 45180
$ws.Cells.Item(348, 5).Value  = 13
$ws.Cells.Item(348, 6).Value  = 100112001
$ws.Cells.Item(348, 7).Value  = "Berenjena"
$ws.Cells.Item(348, 8).Value  = "Sin especificar"
$ws.Cells.Item(348, 9).Value  = "Primera"
$ws.Cells.Item(348, 10).Value = 70
$ws.Cells.Item(348, 11).Value = 9000
$ws.Cells.Item(348, 12).Value = 10000
$ws.Cells.Item(348, 13).Value = 9500
$ws.Cells.Item(348, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(348, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(348, 16).Value = 190
$ws.Cells.Item(348, 17).Value = 50
$ws.Cells.Item(348, 18).Value = "Hortaliza"
